# Generate Report for Handoff
# A new source file (4cf71334-16e8-4d7e-ab7b-58358ed7f1ce.md) has been
# handed off alongside the existing 6403701d-b5be-4da2-bdab-d61aa05cb8b0.md
# file. This inserts a new row (row 2) on every sheet for the new file,
# pushing the previously-existing rows down by one.

$wb = $excel.ActiveWorkbook

$hyperlinkColor = 15570276  # RGB(0x64,0x95,0xED) as an OLE/VBA color value
$xlUnderlineStyleSingle = 2

# Hyperlinks.Add() always (re)stamps its own default "Hyperlink" look on the
# anchor cell, so the custom blue/underline font that matches the rest of
# the workbook has to be (re)applied *after* the hyperlink exists, not
# before -- otherwise Add() clobbers it. Resetting Font.Name first drops
# the theme/scheme binding Add() attaches, so the run collapses back onto
# the sheet's existing hyperlink font instead of minting a near-duplicate.
function Style-AsHyperlink($range) {
    $range.Font.Name = "Calibri"
    $range.Font.Underline = $xlUnderlineStyleSingle
    $range.Font.Color = $hyperlinkColor
}

function Style-AsDateText($range) {
    $range.NumberFormat = "yyyy-mm-dd HH:mm:ss"
}

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# Make room for the new file's row; existing rows 2-3 shift down to 3-4
# (styles/number-formats travel with the shifted cells automatically).
$ws.Rows.Item(2).Insert()

$ws.Range("A2").Value = "4cf71334-16e8-4d7e-ab7b-58358ed7f1ce.md"
$ws.Range("B2").Value = "Ready for handoff"
$ws.Range("C2").Value = "Ready for handoff"

# Rebuild every hyperlink (row-insert does not itself relocate them).
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/e51e32d86c2cf021cc41642040412fbedde74141/e2e/4cf71334-16e8-4d7e-ab7b-58358ed7f1ce.md", [Type]::Missing, [Type]::Missing, "4cf71334-16e8-4d7e-ab7b-58358ed7f1ce.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/e51e32d86c2cf021cc41642040412fbedde74141/e2e/6403701d-b5be-4da2-bdab-d61aa05cb8b0.md", [Type]::Missing, [Type]::Missing, "6403701d-b5be-4da2-bdab-d61aa05cb8b0.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/e51e32d86c2cf021cc41642040412fbedde74141/.localization-config", [Type]::Missing, [Type]::Missing, ".localization-config") | Out-Null

Style-AsHyperlink $ws.Range("A2")
Style-AsHyperlink $ws.Range("A3")
Style-AsHyperlink $ws.Range("A4")

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Rows.Item(2).Insert()

$ws.Range("A2").Value = "4cf71334-16e8-4d7e-ab7b-58358ed7f1ce.md"
$ws.Range("B2").Value = "Ready for handoff"
$ws.Range("C2").Value = "4cf71334-16e8-4d7e-ab7b-58358ed7f1ce.03cee19d3a21500d37a6a1843b4252d410dfe297.zh-cn.xlf"
$ws.Range("D2").Value = "2016-03-11 05:59:06"
$ws.Range("G2").Value = "0001-01-01 00:00:00"
$ws.Range("H2").Value = "Include"
Style-AsDateText $ws.Range("D2")

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/e51e32d86c2cf021cc41642040412fbedde74141/e2e/4cf71334-16e8-4d7e-ab7b-58358ed7f1ce.md", [Type]::Missing, [Type]::Missing, "4cf71334-16e8-4d7e-ab7b-58358ed7f1ce.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/03cee19d3a21500d37a6a1843b4252d410dfe297/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/4cf71334-16e8-4d7e-ab7b-58358ed7f1ce.03cee19d3a21500d37a6a1843b4252d410dfe297.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "4cf71334-16e8-4d7e-ab7b-58358ed7f1ce.03cee19d3a21500d37a6a1843b4252d410dfe297.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/e51e32d86c2cf021cc41642040412fbedde74141/e2e/6403701d-b5be-4da2-bdab-d61aa05cb8b0.md", [Type]::Missing, [Type]::Missing, "6403701d-b5be-4da2-bdab-d61aa05cb8b0.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dd64127422106895c33b1c0f7a1ccab3955b4b9d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/6403701d-b5be-4da2-bdab-d61aa05cb8b0.19bc49e05b5f45a10969064ff5893f4fa0807afd.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "6403701d-b5be-4da2-bdab-d61aa05cb8b0.19bc49e05b5f45a10969064ff5893f4fa0807afd.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/e51e32d86c2cf021cc41642040412fbedde74141/.localization-config", [Type]::Missing, [Type]::Missing, ".localization-config") | Out-Null

Style-AsHyperlink $ws.Range("A2")
Style-AsHyperlink $ws.Range("C2")
Style-AsHyperlink $ws.Range("A3")
Style-AsHyperlink $ws.Range("C3")
Style-AsHyperlink $ws.Range("A4")

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Rows.Item(2).Insert()

$ws.Range("A2").Value = "4cf71334-16e8-4d7e-ab7b-58358ed7f1ce.md"
$ws.Range("B2").Value = "Ready for handoff"
$ws.Range("C2").Value = "4cf71334-16e8-4d7e-ab7b-58358ed7f1ce.03cee19d3a21500d37a6a1843b4252d410dfe297.de-de.xlf"
$ws.Range("D2").Value = "2016-03-11 05:59:14"
$ws.Range("G2").Value = "0001-01-01 00:00:00"
$ws.Range("H2").Value = "Include"
Style-AsDateText $ws.Range("D2")

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/e51e32d86c2cf021cc41642040412fbedde74141/e2e/4cf71334-16e8-4d7e-ab7b-58358ed7f1ce.md", [Type]::Missing, [Type]::Missing, "4cf71334-16e8-4d7e-ab7b-58358ed7f1ce.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/03cee19d3a21500d37a6a1843b4252d410dfe297/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/4cf71334-16e8-4d7e-ab7b-58358ed7f1ce.03cee19d3a21500d37a6a1843b4252d410dfe297.de-de.xlf", [Type]::Missing, [Type]::Missing, "4cf71334-16e8-4d7e-ab7b-58358ed7f1ce.03cee19d3a21500d37a6a1843b4252d410dfe297.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/e51e32d86c2cf021cc41642040412fbedde74141/e2e/6403701d-b5be-4da2-bdab-d61aa05cb8b0.md", [Type]::Missing, [Type]::Missing, "6403701d-b5be-4da2-bdab-d61aa05cb8b0.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5319bec0e0396c4d63fc6abf2fca5f934d671068/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/6403701d-b5be-4da2-bdab-d61aa05cb8b0.19bc49e05b5f45a10969064ff5893f4fa0807afd.de-de.xlf", [Type]::Missing, [Type]::Missing, "6403701d-b5be-4da2-bdab-d61aa05cb8b0.19bc49e05b5f45a10969064ff5893f4fa0807afd.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/e51e32d86c2cf021cc41642040412fbedde74141/.localization-config", [Type]::Missing, [Type]::Missing, ".localization-config") | Out-Null

Style-AsHyperlink $ws.Range("A2")
Style-AsHyperlink $ws.Range("C2")
Style-AsHyperlink $ws.Range("A3")
Style-AsHyperlink $ws.Range("C3")
Style-AsHyperlink $ws.Range("A4")

Write-Output "Done: inserted handoff rows for 4cf71334-16e8-4d7e-ab7b-58358ed7f1ce on all sheets."
